$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# FilesTab query (cell B4): drop the "File Type" and "Breed" output columns
# from the Cypher RETURN clause (ICDC Breed 1-14 script correction).
$newQuery = "MATCH (f:file)-->(parent)`nWITH DISTINCT f, parent`nMATCH (f)-[*]->(c:case)<--(demo:demographic)`nWHERE demo.breed IN ['Dalmatian']  `nOPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`nOPTIONAL MATCH (samp:sample)-->(c)`nWITH DISTINCT f, parent, c, demo, diag, s`nRETURN  coalesce(f.file_name, '') AS ``File Name``,`n        coalesce(labels(parent)[0], '') AS ``Association``,`n        coalesce(f.file_description, '') AS ``Description``,`n        coalesce(f.file_format, '') AS ``Format``,`n        coalesce(f.file_size, '') AS ``Size``,`n        coalesce(c.case_id, '') AS ``Case ID``,`n        coalesce(diag.disease_term,'') AS Diagnosis , `n        coalesce(s.clinical_study_designation,'') AS ``Study Code``"

$ws.Range("B4").Value = $newQuery

# The cell's wrapped text shrank from 17 to 15 lines, so the row autosizes
# from 246.5pt down to 217.5pt (14.5pt/line).
$ws.Rows.Item(4).RowHeight = 217.5

# Selection moves from A4 to B4.
$ws.Range("B4").Select()
